# 100ppi spot number upgrade to 45, and modified accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label change (F1) ---
# shared string order must be preserved to match target file, so string
# cells are written in the same order they first appear in the final
# sharedStrings table: F1, A2, A3, A4, C2, C3, C4
$ws.Range("F1").Value = "顶/底指标（基差*极限）"

# --- Row 2 (existing row, values updated in place) ---
$ws.Range("A2").Value = "镍"

# --- Row 3 (new row) ---
# Copy the formatting (border/bold/alignment) of A2 down to the new code cells
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("A3").Value = "锡"

# --- Row 4 (new row) ---
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A4").Value = "PTA"

# Assign these as a text formula first, then collapse the formula down to a
# literal value; this yields a plain shared-string cell (t="s") without
# nudging NumberFormat/Style (which would otherwise leave a stray style
# behind in styles.xml).
$ws.Range("C2").Formula = "=""1807"""
$ws.Range("C2").Value = $ws.Range("C2").Value

$ws.Range("C3").Formula = "=""1805"""
$ws.Range("C3").Value = $ws.Range("C3").Value

$ws.Range("C4").Formula = "=""1809"""
$ws.Range("C4").Value = $ws.Range("C4").Value

# --- Numeric columns ---
$ws.Range("B2").Value = 2.016546018614271
$ws.Range("D2").Value = 98920
$ws.Range("E2").Value = 0.1428571428571428
$ws.Range("F2").Value = -0.03909926869886322

$ws.Range("B3").Value = 0.351493848857645
$ws.Range("D3").Value = 142890
$ws.Range("E3").Value = 0.05185185185185185
$ws.Range("F3").Value = -0.005085295453071285

$ws.Range("B4").Value = 0.5599637681159394
$ws.Range("D4").Value = 5432
$ws.Range("E4").Value = 0.2668297280315947
$ws.Range("F4").Value = 0.5715949810073824
